$d = $word.ActiveDocument

# Helper: find an exact run of text and overwrite it via Range.Text
# (avoids Word's "smart quotes" autocorrect that Find.Execute's Replace
#  With argument would otherwise trigger on straight apostrophes).
function Replace-PlainText($old, $new) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Text = $old
    $rng.Find.MatchCase = $true
    $rng.Find.MatchWholeWord = $false
    $rng.Find.MatchWildcards = $false
    $rng.Find.Forward = $true
    $rng.Find.Wrap = 1
    if ($rng.Find.Execute()) {
        $rng.Text = $new
    }
}

# --- Header block (right-aligned contact info) ---
Replace-PlainText "Guillermo x" "Guillermo Garcia"
Replace-PlainText "x, x" "Valencia, Spain"

# The two remaining "x" placeholders (email, phone) are standalone runs
# with identical text, so find+replace them in document order.
$rngEmail = $d.Content
$rngEmail.Find.ClearFormatting()
$rngEmail.Find.Text = "x"
$rngEmail.Find.MatchWholeWord = $true
$rngEmail.Find.MatchCase = $true
$rngEmail.Find.MatchWildcards = $false
$rngEmail.Find.Forward = $true
$rngEmail.Find.Wrap = 1
if ($rngEmail.Find.Execute()) {
    $rngEmail.Text = "guillermo@test.com"
}

$rngPhone = $d.Content
$rngPhone.Find.ClearFormatting()
$rngPhone.Find.Text = "x"
$rngPhone.Find.MatchWholeWord = $true
$rngPhone.Find.MatchCase = $true
$rngPhone.Find.MatchWildcards = $false
$rngPhone.Find.Forward = $true
$rngPhone.Find.Wrap = 1
if ($rngPhone.Find.Execute()) {
    $rngPhone.Text = "1234567899"
}

# --- Body paragraph replacements (first three: no new breaks needed) ---
Replace-PlainText `
    "I am writing to express my interest in the Software Engineer position at Swish Analytics. As a highly experienced developer with a Bachelor's degree in Computer Science, I believe I possess the skills and expertise necessary to contribute effectively to your Visual Services team." `
    "I am writing to express my interest in the Java Developer position at B.E.A.T. LLC, as advertised. With a Bachelor's degree in Digital Media Software Engineering and hands-on experience in Java, JavaScript, HTML, and CSS, I am excited about the opportunity to contribute to your team."

Replace-PlainText `
    "In my previous roles, I have successfully designed and developed data analytics platforms utilizing cutting-edge technologies. I have a solid background in backend Restful API development, with experience in frameworks like Django and Flask. Additionally, my expertise in front-end development using JavaScript, ReactJS, and NodeJS aligns perfectly with the requirements of the position at Swish Analytics." `
    "In my role as a Software Engineer Intern Consultant at Masetto Logistics, I designed and developed features in Python for real-time fleet management, incorporating API calls for precise location tracking. I also collaborated with cross-functional teams in daily standup meetings and biweekly scrum gatherings, demonstrating my ability to work effectively in an agile environment."

Replace-PlainText `
    "I am excited about the opportunity to work closely with key stakeholders and product owners to drive technical design for various use cases. My ability to drive projects end-to-end and produce high-quality software that is rigorously tested and reviewed makes me a strong candidate for this role. I am confident that my experience in SQL & data analytics, including time series analytical queries and data modeling, will enable me to excel in this position." `
    "My coursework in Machine Learning, Java, React, UI Design, and Agile methodologies has equipped me with the skills necessary to excel in this position. Additionally, my experience collaborating with DevOps teams, conducting functional tests, and leveraging Jira for test case management aligns with the qualifications for this role."

# --- Final paragraph: replace the text AND append a brand-new sentence
#     (with its own blank line) ahead of the pre-existing trailing blank
#     line. Using the Find "Replace With" special code ^l (manual line
#     break) keeps the newly-inserted <w:br/><w:br/> as their own run,
#     separate from the original trailing <w:br/><w:br/>.
$rngFinal = $d.Content
$rngFinal.Find.ClearFormatting()
$oldFinal = "I am eager to bring my innovative approach and technical skills to Swish Analytics to help build the next-generation data analytics platform. I am confident that my background and expertise make me a perfect fit for this role. Thank you for considering my application. I look forward to the opportunity to discuss how my qualifications align with the needs of your team."
$newFinal = "I am confident that my technical skills, experience, and strong communication abilities make me a qualified candidate for the Java Developer position at B.E.A.T. LLC. I am eager to bring my expertise in Java, HTML, and JavaScript to your team and contribute to the development of complex business solutions.^l^lThank you for considering my application. I look forward to the possibility of discussing how my background, skills, and enthusiasm can contribute to the success of your team."
$rngFinal.Find.Execute($oldFinal, $true, $true, $false, $false, $false, $true, 1, $false, $newFinal, 2) | Out-Null
